$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.406.55'
$ws.Range("E2").Value = '  -7.99%  '

$ws.Range("D3").Value = '1.677.63'
$ws.Range("E3").Value = '  -6.89%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.36%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -6.48%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.006'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.28%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4961'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -16.68%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2591'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.67%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '21.67'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -7.04%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06145'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -9.93%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07292'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.29%  '

$ws.Range("D12").Value = '1.728.60'
$ws.Range("E12").Value = '  -4.06%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.408'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -7.32%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5724'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -8.20%  '

$ws.Range("D15").Value = '1.904.46'
$ws.Range("E15").Value = '  -6.96%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008157'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -12.19%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.25'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -14.74%  '

$ws.Range("D18").Value = '26.427.87'
$ws.Range("E18").Value = '  -7.77%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.977'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -8.99%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.006'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.30%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.71'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.31%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '183.59'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -12.53%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.164'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -9.92%  '

$ws.Range("E24").Value = '  +0.29%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.74'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.87%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.479'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.64%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1127'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -11.42%  '

$ws.Range("E28").Value = '  -5.05%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.302'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -8.76%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05684'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.47%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.315'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -7.39%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.470'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -8.14%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.455'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -7.72%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.627'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.22%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.000'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.84%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.364'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.89%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5870'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -7.94%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.627'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.11%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01581'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.50%  '

$ws.Range("D40").Value = '1.068.69'
$ws.Range("E40").Value = '  -5.50%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.871'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -8.91%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8482'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.51%  '

$ws.Range("E43").Value = '  -0.30%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '97.72'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.90%  '

$ws.Range("D45").Value = '1.834.17'
$ws.Range("E45").Value = '  -6.44%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.15'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.14%  '

$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.005'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.14%  '

$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000104'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.80%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.023'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.88%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4322'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.74%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05178'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.35%  '
